{"js": "// Update the division-problem table: replace each cell's old expression\n// with its new expression, addressed by (row, col) position so that\n// duplicate values (e.g. \"912\u00f76=\" appears both before and after editing)\n// never get confused with one another.\nconst table = context.document.body.tables.getFirst();\n\n// [rowIndex, columnIndex, expectedOldText, newText]\nconst updates = [\n  [0, 0, \"335\u00f79=\", \"332\u00f73=\"],\n  [0, 1, \"833\u00f79=\", \"707\u00f78=\"],\n  [0, 2, \"700\u00f75=\", \"998\u00f73=\"],\n  [0, 3, \"724\u00f76=\", \"118\u00f78=\"],\n  [0, 4, \"961\u00f72=\", \"640\u00f76=\"],\n  [4, 0, \"808\u00f75=\", \"937\u00f79=\"],\n  [4, 1, \"932\u00f78=\", \"366\u00f72=\"],\n  [4, 2, \"933\u00f79=\", \"908\u00f74=\"],\n  [4, 3, \"384\u00f73=\", \"109\u00f72=\"],\n  [4, 4, \"759\u00f77=\", \"920\u00f75=\"],\n  [8, 0, \"754\u00f75=\", \"132\u00f76=\"],\n  [8, 1, \"131\u00f73=\", \"198\u00f75=\"],\n  [8, 2, \"871\u00f79=\", \"646\u00f76=\"],\n  [8, 3, \"981\u00f77=\", \"681\u00f77=\"],\n  [8, 4, \"846\u00f78=\", \"601\u00f74=\"],\n  [12, 0, \"647\u00f77=\", \"622\u00f79=\"],\n  [12, 1, \"830\u00f75=\", \"706\u00f75=\"],\n  [12, 2, \"629\u00f77=\", \"765\u00f73=\"],\n  [12, 3, \"588\u00f72=\", \"254\u00f79=\"],\n  [12, 4, \"690\u00f76=\", \"912\u00f76=\"],\n  [16, 0, \"994\u00f79=\", \"103\u00f72=\"],\n  [16, 1, \"912\u00f76=\", \"310\u00f74=\"],\n  [16, 2, \"123\u00f73=\", \"308\u00f74=\"],\n  [16, 3, \"282\u00f76=\", \"931\u00f75=\"],\n  [16, 4, \"620\u00f79=\", \"214\u00f72=\"],\n];\n\n// Grab each target paragraph range up front.\nconst ranges = updates.map(([row, col]) => {\n  const cell = table.getCell(row, col);\n  return cell.body.paragraphs.getFirst().getRange();\n});\nranges.forEach((range) => range.load(\"text\"));\nawait context.sync();\n\n// Verify we are about to overwrite the expected old value, then replace it.\nranges.forEach((range, i) => {\n  const [, , oldText, newText] = updates[i];\n  if (range.text !== oldText) {\n    throw new Error(\n      `Unexpected cell text \"${range.text}\"; expected \"${oldText}\".`\n    );\n  }\n  range.insertText(newText, Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Update the division-problem table: replace each cell's old expression\n# with its new expression, addressed by (row, col) position so that\n# duplicate values (e.g. \"912\u00f76=\" appears both before and after editing)\n# never get confused with one another.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each entry: row (1-based), column (1-based), expected old text, new text\n$updates = @(\n  @(1, 1, \"335\u00f79=\", \"332\u00f73=\"),\n  @(1, 2, \"833\u00f79=\", \"707\u00f78=\"),\n  @(1, 3, \"700\u00f75=\", \"998\u00f73=\"),\n  @(1, 4, \"724\u00f76=\", \"118\u00f78=\"),\n  @(1, 5, \"961\u00f72=\", \"640\u00f76=\"),\n  @(5, 1, \"808\u00f75=\", \"937\u00f79=\"),\n  @(5, 2, \"932\u00f78=\", \"366\u00f72=\"),\n  @(5, 3, \"933\u00f79=\", \"908\u00f74=\"),\n  @(5, 4, \"384\u00f73=\", \"109\u00f72=\"),\n  @(5, 5, \"759\u00f77=\", \"920\u00f75=\"),\n  @(9, 1, \"754\u00f75=\", \"132\u00f76=\"),\n  @(9, 2, \"131\u00f73=\", \"198\u00f75=\"),\n  @(9, 3, \"871\u00f79=\", \"646\u00f76=\"),\n  @(9, 4, \"981\u00f77=\", \"681\u00f77=\"),\n  @(9, 5, \"846\u00f78=\", \"601\u00f74=\"),\n  @(13, 1, \"647\u00f77=\", \"622\u00f79=\"),\n  @(13, 2, \"830\u00f75=\", \"706\u00f75=\"),\n  @(13, 3, \"629\u00f77=\", \"765\u00f73=\"),\n  @(13, 4, \"588\u00f72=\", \"254\u00f79=\"),\n  @(13, 5, \"690\u00f76=\", \"912\u00f76=\"),\n  @(17, 1, \"994\u00f79=\", \"103\u00f72=\"),\n  @(17, 2, \"912\u00f76=\", \"310\u00f74=\"),\n  @(17, 3, \"123\u00f73=\", \"308\u00f74=\"),\n  @(17, 4, \"282\u00f76=\", \"931\u00f75=\"),\n  @(17, 5, \"620\u00f79=\", \"214\u00f72=\"),\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $oldText = $u[2]\n    $newText = $u[3]\n    $cell = $t.Cell($row, $col)\n    # Cell text carries a trailing cell-mark (CR + BEL); strip it before comparing.\n    $currentText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($currentText -ne $oldText) {\n        throw \"Unexpected cell text at row ${row}, col ${col}: expected `\"$oldText`\" but found `\"$currentText`\".\"\n    }\n    $cell.Range.Text = $newText\n}\n"}
